$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4362
$ws1.Range("F3").Value = 2454
$ws1.Range("F6").Value = 49
$ws1.Range("F10").Value = 150
$ws1.Range("F11").Value = 160
$ws1.Range("F12").Value = 1627
$ws1.Range("F14").Value = 3439

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4362
$ws4.Range("F3").Value = 2454
$ws4.Range("F7").Value = 49
$ws4.Range("F12").Value = 150
$ws4.Range("F13").Value = 160
$ws4.Range("F16").Value = 1627
$ws4.Range("F18").Value = 3439
